$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 507.08334
$ws.Range("I19").Value = 445.5
$ws.Range("J19").Value = 593.3
$ws.Range("K19").Value = 445.5
$ws.Range("L19").Value = 593.3
$ws.Range("M19").Value = -270.5
$ws.Range("N19").Value = -943.3
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null
$ws.Range("H41").Value = 954.75
$ws.Range("I41").Value = 567.8
$ws.Range("K41").Value = 567.8
$ws.Range("M41").Value = -127.8
$ws.Range("H64").Value = 3307.9534
$ws.Range("I64").Value = 2732.5173
$ws.Range("J64").Value = 4499.9287
$ws.Range("K64").Value = 2732.5173
$ws.Range("L64").Value = 4499.9287
$ws.Range("M64").Value = -2484.5173
$ws.Range("N64").Value = -4995.9287
$ws.Range("H67").Value = 3307.9534
$ws.Range("I67").Value = 2732.5173
$ws.Range("J67").Value = 4499.9287
$ws.Range("K67").Value = 2732.5173
$ws.Range("L67").Value = 4499.9287
$ws.Range("M67").Value = -1874.5173
$ws.Range("N67").Value = -6215.9287
$ws.Range("H106").Value = 1794.1765
$ws.Range("I106").Value = 1818.8125
$ws.Range("K106").Value = 1818.8125
$ws.Range("M106").Value = -1187.8125
$ws.Range("H111").Value = 792.4
$ws.Range("I111").Value = 1003.3333
$ws.Range("J111").Value = 476
$ws.Range("K111").Value = 3009.9999
$ws.Range("L111").Value = 1428
$ws.Range("M111").Value = 57.0001000000002
$ws.Range("N111").Value = -7562
$ws.Range("H112").Value = 2443.7273
$ws.Range("J112").Value = 2498
$ws.Range("L112").Value = 7494
$ws.Range("N112").Value = -9710
$ws.Range("H135").Value = 4348
$ws.Range("I135").Value = 3386.6667
$ws.Range("J135").Value = 13000
$ws.Range("K135").Value = 30480.0003
$ws.Range("L135").Value = 117000
$ws.Range("M135").Value = -27945.0003
$ws.Range("N135").Value = -122070
$ws.Range("H137").Value = 12005.412
$ws.Range("I137").Value = 2637
$ws.Range("J137").Value = 20332.889
$ws.Range("K137").Value = 7911
$ws.Range("L137").Value = 60998.667
$ws.Range("M137").Value = -5361
$ws.Range("N137").Value = -66098.667
$ws.Range("H138").Value = 2293.1968
$ws.Range("J138").Value = 2579.0613
$ws.Range("L138").Value = 7737.1839
$ws.Range("N138").Value = -18017.1839

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H26").Value = 340166.66
$ws.Range("I26").Value = 340166.66
$ws.Range("K26").Value = 340166.66
$ws.Range("M26").Value = -339836.66
$ws.Range("H32").Value = 142942.12
$ws.Range("I32").Value = 151982.1
$ws.Range("J32").Value = 19998.6
$ws.Range("K32").Value = 151982.1
$ws.Range("L32").Value = 19998.6
$ws.Range("M32").Value = -151695.1
$ws.Range("N32").Value = -20572.6
$ws.Range("H45").Value = 2007.3334
$ws.Range("I45").Value = 2011
$ws.Range("K45").Value = 2011
$ws.Range("M45").Value = -1634
$ws.Range("H61").Value = 9114
$ws.Range("I61").Value = 10393.875
$ws.Range("K61").Value = 10393.875
$ws.Range("M61").Value = -10181.875
$ws.Range("H63").Value = 3538.077
$ws.Range("I63").Value = 3090.4546
$ws.Range("K63").Value = 3090.4546
$ws.Range("M63").Value = -2404.4546
$ws.Range("H66").Value = 3538.077
$ws.Range("I66").Value = 3090.4546
$ws.Range("K66").Value = 15452.273
$ws.Range("M66").Value = -12020.273
$ws.Range("H74").Value = 13935.833
$ws.Range("I74").Value = 1761.6666
$ws.Range("J74").Value = 26110
$ws.Range("K74").Value = 1761.6666
$ws.Range("L74").Value = 26110
$ws.Range("M74").Value = -887.6666
$ws.Range("N74").Value = -27858
$ws.Range("H77").Value = 13935.833
$ws.Range("I77").Value = 1761.6666
$ws.Range("J77").Value = 26110
$ws.Range("K77").Value = 8808.333000000001
$ws.Range("L77").Value = 130550
$ws.Range("M77").Value = -4440.333000000001
$ws.Range("N77").Value = -139286
$ws.Range("H88").Value = 1738.4445
$ws.Range("J88").Value = 1942.3889
$ws.Range("L88").Value = 1942.3889
$ws.Range("N88").Value = -2754.3889
$ws.Range("H91").Value = 1738.4445
$ws.Range("J91").Value = 1942.3889
$ws.Range("L91").Value = 1942.3889
$ws.Range("N91").Value = -4750.3889
$ws.Range("H110").Value = 1432.1111
$ws.Range("I110").Value = 1148.1666
$ws.Range("K110").Value = 1148.1666
$ws.Range("M110").Value = 896.8334
$ws.Range("H122").Value = 33334584
$ws.Range("I122").Value = 38462480
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 115387440
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -115384990
$ws.Range("N122").Value = -14650
$ws.Range("H125").Value = 66999.5
$ws.Range("J125").Value = 66999.5
$ws.Range("L125").Value = 66999.5
$ws.Range("N125").Value = -76839.5
$ws.Range("H132").Value = 1138358.6
$ws.Range("I132").Value = 1317665.5
$ws.Range("J132").Value = 2748.3333
$ws.Range("K132").Value = 3952996.5
$ws.Range("L132").Value = 8244.999899999999
$ws.Range("M132").Value = -3950466.5
$ws.Range("N132").Value = -13304.9999
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
$ws.Range("H136").Value = 9114
$ws.Range("I136").Value = 10393.875
$ws.Range("K136").Value = 31181.625
$ws.Range("M136").Value = -28631.625

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47533.523
$ws.Range("I20").Value = 77135.78999999999
$ws.Range("J20").Value = 1485.5555
$ws.Range("K20").Value = 77135.78999999999
$ws.Range("L20").Value = 1485.5555
$ws.Range("M20").Value = -76888.78999999999
$ws.Range("N20").Value = -1979.5555
$ws.Range("H50").Value = 440000
$ws.Range("J50").Value = 440000
$ws.Range("L50").Value = 440000
$ws.Range("N50").Value = -441148
$ws.Range("H86").Value = 2305.5833
$ws.Range("I86").Value = 2475.125
$ws.Range("J86").Value = 1966.5
$ws.Range("K86").Value = 2475.125
$ws.Range("L86").Value = 1966.5
$ws.Range("M86").Value = -1352.125
$ws.Range("N86").Value = -4212.5
$ws.Range("H89").Value = 2305.5833
$ws.Range("I89").Value = 2475.125
$ws.Range("J89").Value = 1966.5
$ws.Range("K89").Value = 12375.625
$ws.Range("L89").Value = 9832.5
$ws.Range("M89").Value = -6759.625
$ws.Range("N89").Value = -21064.5
$ws.Range("H94").Value = 3792.6
$ws.Range("I94").Value = 3792.6
$ws.Range("K94").Value = 3792.6
$ws.Range("M94").Value = -3341.6
$ws.Range("H99").Value = 8223.143
$ws.Range("I99").Value = 17325.666
$ws.Range("J99").Value = 1396.25
$ws.Range("K99").Value = 17325.666
$ws.Range("L99").Value = 1396.25
$ws.Range("M99").Value = -15827.666
$ws.Range("N99").Value = -4392.25
$ws.Range("H105").Value = 3171.5952
$ws.Range("I105").Value = 2474.4517
$ws.Range("K105").Value = 2474.4517
$ws.Range("M105").Value = -727.4517000000001
$ws.Range("H107").Value = 16667681
$ws.Range("I107").Value = 17858160
$ws.Range("J107").Value = 979
$ws.Range("K107").Value = 17858160
$ws.Range("L107").Value = 979
$ws.Range("M107").Value = -17856240
$ws.Range("N107").Value = -4819
$ws.Range("H134").Value = 11283.15
$ws.Range("I134").Value = 6613.8423
$ws.Range("K134").Value = 19841.5269
$ws.Range("M134").Value = -17306.5269

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1474.3334
$ws.Range("I21").Value = 1474.3334
$ws.Range("K21").Value = 1474.3334
$ws.Range("M21").Value = -1239.3334
$ws.Range("H58").Value = 12607.357
$ws.Range("I58").Value = 5601.5557
$ws.Range("K58").Value = 5601.5557
$ws.Range("M58").Value = -5398.5557
$ws.Range("H99").Value = 19311.924
$ws.Range("I99").Value = 22096
$ws.Range("J99").Value = 3999.5
$ws.Range("K99").Value = 22096
$ws.Range("L99").Value = 3999.5
$ws.Range("M99").Value = -20598
$ws.Range("N99").Value = -6995.5
$ws.Range("H126").Value = 19311.924
$ws.Range("I126").Value = 22096
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 66288
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -63818
$ws.Range("N126").Value = -16938.5
$ws.Range("H132").Value = 3317.5334
$ws.Range("I132").Value = 3368.7856
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 10106.3568
$ws.Range("L132").Value = 7800
$ws.Range("M132").Value = -7576.356800000001
$ws.Range("N132").Value = -12860
$ws.Range("H134").Value = 3391.0625
$ws.Range("I134").Value = 2950.4666
$ws.Range("K134").Value = 8851.399800000001
$ws.Range("M134").Value = -6316.399800000001
$ws.Range("H136").Value = 12607.357
$ws.Range("I136").Value = 5601.5557
$ws.Range("K136").Value = 16804.6671
$ws.Range("M136").Value = -14254.6671

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1864.6666
$ws.Range("J5").Value = 2379.75
$ws.Range("L5").Value = 7139.25
$ws.Range("N5").Value = -7363.25
$ws.Range("H14").Value = 357872.5
$ws.Range("I14").Value = 357872.5
$ws.Range("K14").Value = 1073617.5
$ws.Range("M14").Value = -1073444.5
$ws.Range("H68").Value = 9799
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 9799
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 29397
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -31019
$ws.Range("H71").Value = 9799
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 9799
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 88191
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -96303
$ws.Range("H119").Value = 11184.454
$ws.Range("I119").Value = 3605.8
$ws.Range("K119").Value = 10817.4
$ws.Range("M119").Value = -5979.400000000001
$ws.Range("H120").Value = 17537.066
$ws.Range("I120").Value = 12305.6
$ws.Range("K120").Value = 36916.8
$ws.Range("M120").Value = -32078.8
$ws.Range("H121").Value = 3575438.5
$ws.Range("I121").Value = 808.1429000000001
$ws.Range("K121").Value = 2424.4287
$ws.Range("M121").Value = -1114.4287
$ws.Range("H133").Value = 10300.682
$ws.Range("I133").Value = 8107.6665
$ws.Range("K133").Value = 24322.9995
$ws.Range("M133").Value = -19262.9995
$ws.Range("H134").Value = 7704.7896
$ws.Range("I134").Value = 2139.1
$ws.Range("K134").Value = 6417.299999999999
$ws.Range("M134").Value = -1347.299999999999
$ws.Range("H135").Value = 1864.6666
$ws.Range("J135").Value = 2379.75
$ws.Range("L135").Value = 21417.75
$ws.Range("N135").Value = -26487.75
$ws.Range("H139").Value = 6145.048
$ws.Range("I139").Value = 2087.1667
$ws.Range("J139").Value = 11555.556
$ws.Range("K139").Value = 6261.500100000001
$ws.Range("L139").Value = 34666.66800000001
$ws.Range("M139").Value = -1121.500100000001
$ws.Range("N139").Value = -44946.66800000001
$ws.Range("H141").Value = 15425.23
$ws.Range("I141").Value = 12007
$ws.Range("J141").Value = 16944.445
$ws.Range("K141").Value = 36021
$ws.Range("L141").Value = 50833.335
$ws.Range("M141").Value = -30841
$ws.Range("N141").Value = -61193.335

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10889.667
$ws.Range("I12").Value = 14999
$ws.Range("J12").Value = 8835
$ws.Range("K12").Value = 14999
$ws.Range("L12").Value = 8835
$ws.Range("M12").Value = -14859
$ws.Range("N12").Value = -9115
$ws.Range("H97").Value = 2507.7407
$ws.Range("I97").Value = 2060.5
$ws.Range("K97").Value = 2060.5
$ws.Range("M97").Value = -1564.5
$ws.Range("H132").Value = 7077.136
$ws.Range("J132").Value = 17563.143
$ws.Range("L132").Value = 52689.429
$ws.Range("N132").Value = -57749.429

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4576.0713
$ws.Range("I7").Value = 4539.091
$ws.Range("J7").Value = 4711.6665
$ws.Range("K7").Value = 4539.091
$ws.Range("L7").Value = 4711.6665
$ws.Range("M7").Value = -4427.091
$ws.Range("N7").Value = -4935.6665
$ws.Range("H16").Value = 902.82355
$ws.Range("I16").Value = 746
$ws.Range("J16").Value = 1634.6666
$ws.Range("K16").Value = 746
$ws.Range("L16").Value = 1634.6666
$ws.Range("M16").Value = -576
$ws.Range("N16").Value = -1974.6666
$ws.Range("H93").Value = 2547.0833
$ws.Range("I93").Value = 2187.5908
$ws.Range("K93").Value = 2187.5908
$ws.Range("M93").Value = -939.5907999999999
$ws.Range("H122").Value = 3190.3845
$ws.Range("I122").Value = 3094
$ws.Range("J122").Value = 3250.625
$ws.Range("K122").Value = 9282
$ws.Range("L122").Value = 9751.875
$ws.Range("M122").Value = -6832
$ws.Range("N122").Value = -14651.875
$ws.Range("H126").Value = 4576.0713
$ws.Range("I126").Value = 4539.091
$ws.Range("J126").Value = 4711.6665
$ws.Range("K126").Value = 13617.273
$ws.Range("L126").Value = 14134.9995
$ws.Range("M126").Value = -11147.273
$ws.Range("N126").Value = -19074.9995
$ws.Range("H132").Value = 8367368
$ws.Range("I132").Value = 25094444
$ws.Range("J132").Value = 3830
$ws.Range("K132").Value = 75283332
$ws.Range("L132").Value = 11490
$ws.Range("M132").Value = -75280802
$ws.Range("N132").Value = -16550
$ws.Range("H133").Value = 59913.5
$ws.Range("J133").Value = 59913.5
$ws.Range("L133").Value = 59913.5
$ws.Range("N133").Value = -64973.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 45550
$ws.Range("J92").Value = 45550
$ws.Range("L92").Value = 45550
$ws.Range("N92").Value = -50542
$ws.Range("H107").Value = 2631.8
$ws.Range("I107").Value = 2254.238
$ws.Range("J107").Value = 3512.7778
$ws.Range("K107").Value = 6762.714
$ws.Range("L107").Value = 10538.3334
$ws.Range("M107").Value = -4842.714
$ws.Range("N107").Value = -14378.3334
$ws.Range("H119").Value = 69997.5
$ws.Range("J119").Value = 69997.5
$ws.Range("L119").Value = 69997.5
$ws.Range("N119").Value = -79673.5
$ws.Range("H120").Value = 400000
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = $null
$ws.Range("H122").Value = 189302.67
$ws.Range("H126").Value = 3329.5
$ws.Range("I126").Value = 3331.0833
$ws.Range("J126").Value = 3324.75
$ws.Range("K126").Value = 9993.249899999999
$ws.Range("L126").Value = 9974.25
$ws.Range("M126").Value = -7523.249899999999
$ws.Range("N126").Value = -14914.25
$ws.Range("H132").Value = 3174.7368
$ws.Range("I132").Value = 3002.6667
$ws.Range("J132").Value = 3820
$ws.Range("K132").Value = 9008.000100000001
$ws.Range("L132").Value = 11460
$ws.Range("M132").Value = -6478.000100000001
$ws.Range("N132").Value = -16520
$ws.Range("H136").Value = 1168.0625
$ws.Range("I136").Value = 936.25
$ws.Range("J136").Value = 1863.5
$ws.Range("K136").Value = 2808.75
$ws.Range("L136").Value = 5590.5
$ws.Range("M136").Value = -258.75
$ws.Range("N136").Value = -10690.5
